# Accounts.xlsx — add rows for signed-in accounts table
# Row 3: "12" repeated across A:D (stored as text via quote-prefix, like the
#        existing "1"/"10" header rows)
# Row 4: "1" repeated across A:D (reuses the existing shared string)
# Row 5: the actual account record — Luca / molinari / Username / password

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text storage, matching the existing "1"/"10"
# text cells in rows 1-2 (so these numeric-looking codes stay text, not
# numbers).
$ws.Range("A3:D3").Value = "'12"
$ws.Range("A4:D4").Value = "'1"

# New account record
$ws.Range("A5").Value = "Luca"
$ws.Range("B5").Value = "molinari"
$ws.Range("C5").Value = "Username"
$ws.Range("D5").Value = "password"
